# Daily attendance processing - 2026-01-22 04:52:08
# Normalizes the ordering of names/emails in the "Recorded By" column (G)
# so that the "System"/"system" and previously-first recorder name are
# reordered to match the latest processing pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old "Recorded By" combined strings to their corrected order.
# Any value not present in this map (e.g. single-author entries) is left
# untouched.
$recordedByMap = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

# Determine the last used row from the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $current = $cell.Value2

    if ($null -ne $current -and $recordedByMap.ContainsKey($current)) {
        $cell.Value2 = $recordedByMap[$current]
    }
}
